$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

$ws1.Range("F3").Value = 728
$ws1.Range("F6").Value = 2975
$ws1.Range("F7").Value = 1692
$ws1.Range("F8").Value = 1993
$ws1.Range("F9").Value = 319
$ws1.Range("F10").Value = 298
$ws1.Range("F11").Value = 850
$ws1.Range("F12").Value = 947
$ws1.Range("F14").Value = 418
$ws1.Range("F15").Value = 1135
$ws1.Range("F17").Value = 64
$ws1.Range("F18").Value = 533
$ws1.Range("F19").Value = 7199
$ws1.Range("F21").Value = 1903
$ws1.Range("F23").Value = 187
$ws1.Range("F25").Value = 445
$ws1.Range("G25").Value = "已售罄"
$ws1.Range("F26").Value = 445
$ws1.Range("F28").Value = 1119
$ws1.Range("F29").Value = 943
$ws1.Range("F31").Value = 118
$ws1.Range("F33").Value = 1115
$ws1.Range("F34").Value = 1916
$ws1.Range("F35").Value = 471
$ws1.Range("F36").Value = 9
$ws1.Range("F38").Value = 251
$ws1.Range("F40").Value = 148
$ws1.Range("F41").Value = 270
$ws1.Range("F43").Value = 191
$ws4.Range("F3").Value = 728
$ws4.Range("F9").Value = 2976
$ws4.Range("F10").Value = 1692
$ws4.Range("F11").Value = 1993
$ws4.Range("F12").Value = 319
$ws4.Range("F13").Value = 298
$ws4.Range("F14").Value = 851
$ws4.Range("F16").Value = 947
$ws4.Range("F18").Value = 418
$ws4.Range("F19").Value = 1135
$ws4.Range("F21").Value = 64
$ws4.Range("F22").Value = 533
$ws4.Range("F23").Value = 7199
$ws4.Range("F25").Value = 1905
$ws4.Range("F28").Value = 187
$ws4.Range("F30").Value = 445
$ws4.Range("G30").Value = "已售罄"
$ws4.Range("F31").Value = 445
$ws4.Range("F33").Value = 1119
$ws4.Range("F34").Value = 943
$ws4.Range("F36").Value = 118
$ws4.Range("F37").Value = 1115
$ws4.Range("F38").Value = 1916
$ws4.Range("F39").Value = 471
$ws4.Range("F40").Value = 9
$ws4.Range("F42").Value = 251
$ws4.Range("F44").Value = 148
$ws4.Range("F45").Value = 270
$ws4.Range("F49").Value = 191
